$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to stay as plain text (matches source which stores
# the date as a shared string, not a numeric date serial).
$ws.Range("A2:A4").NumberFormat = "@"

# Row 2 -> GLD / StreetTRACKS Gold Shares
$ws.Range("A2").Value = "2025-12-03"
$ws.Range("B2").Value = "StreetTRACKS Gold Shares"
$ws.Range("C2").Value = "GLD"
$ws.Range("D2").Value = 385.92
$ws.Range("E2").Value = 56.9
$ws.Range("F2").Value = 1.5
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 76
$ws.Range("K2").Value = 66.8
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 65.32892478746797
$ws.Range("O2").Value = "🟢 상승 우위 (다소 완화)"

# Row 3 -> NEM / Newmont Corporation
$ws.Range("A3").Value = "2025-12-03"
$ws.Range("B3").Value = "Newmont Corporation"
$ws.Range("C3").Value = "NEM"
$ws.Range("D3").Value = 89.59999999999999
$ws.Range("E3").Value = 49.8
$ws.Range("F3").Value = 3.85
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 70
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 65.59999999999999
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N3").Value = 65.32892478746797
$ws.Range("O3").Value = "🟢 상승 우위 (다소 완화)"

# Row 4 -> GC=F / Gold Dec 25
$ws.Range("A4").Value = "2025-12-03"
$ws.Range("B4").Value = "Gold Dec 25"
$ws.Range("C4").Value = "GC=F"
$ws.Range("D4").Value = 4226.9
$ws.Range("E4").Value = 71
$ws.Range("F4").Value = 4.08
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 63
$ws.Range("J4").Value = 66
$ws.Range("K4").Value = 56.8
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 65.32892478746797
$ws.Range("O4").Value = "🟢 상승 우위 (다소 완화)"

# Restore the original (unstyled) look for the date cells so no stray
# "text" number-format style is left behind in styles.xml.
$ws.Range("A2:A4").Style = "Normal"
